$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before the current row 5 ("beginDateTime"),
# pushing the existing rows 5-14 down to 7-16.
$ws.Rows("5:6").Insert()

# Populate the two newly inserted rows with the new activity fields.
$ws.Range("A5").Value = "activityName"
$ws.Range("C5").Value = "活动名称"
$ws.Range("A6").Value = "activityDesc"
$ws.Range("C6").Value = "活动描述"

# Update the selected range to match the new data extent.
$ws.Range("C3:C16").Select() | Out-Null

# Add an explicit page setup (paper size / orientation) as in the target file.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
